$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Diamond / value updates across the level grid ---
$ws.Range("AA2").Value = 1
$ws.Range("L6").Value = 3
$ws.Range("C10").Value = 4
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("N12").Value = 3
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("U13").Value = 0
$ws.Range("V13").Value = 0
$ws.Range("Y13").Value = 0
$ws.Range("J14").Value = 3
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 2
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 2
$ws.Range("O15").Value = 2
$ws.Range("P15").Value = 2
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 2
$ws.Range("U15").Value = 2
$ws.Range("V15").Value = 2
$ws.Range("Y15").Value = 2
$ws.Range("Z16").Value = 0
$ws.Range("AA16").Value = 0
$ws.Range("Z17").Value = 0
$ws.Range("AA17").Value = 0
$ws.Range("U18").Value = 4
$ws.Range("S23").Value = 3
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("C26").Value = 1

# --- Added spawn locations: update selection + zoom on the active sheet view ---
$ws.Activate()
$ws.Range("O21").Select()
$excel.ActiveWindow.Zoom = 145
